# cambio de fracciones e historico
# Update quarterly "Personal contratado por honorarios" report: the Q2 2022
# row (with two hired persons) is replaced by a Q3 2022 row stating that no
# personnel was hired by honorarios in that period; the old per-person data
# row is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- 1. Remove the second data row (row 9) -------------------------------
# The sheet used to report two people; now only a single "no hires" row
# remains, so the former row 9 disappears and row 8 becomes the only
# data row.
$ws.Rows("9:9").Delete()

# --- 2. Normalize the style of the fields that become blank --------------
# Columns D,E,F,G,H,I,J,K,L,M,N,O,P,Q,R on row 8 lose their specialised
# formatting (wrap text, shaded fill, hyperlink font, currency format) and
# become plain bordered/left-aligned cells, matching the look of F8 in the
# original sheet.
$ws.Range("F8").Copy()
$ws.Range("D8:R8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Update the period covered by the report ---------------------------
$ws.Range("B8").Value = 44743   # 2022-07-01
$ws.Range("C8").Value = 44834   # 2022-09-30

# --- 4. Clear the per-person / per-contract fields ------------------------
$ws.Range("D8:Q8").ClearContents()

# --- 5. Fill in the note, responsible area and validation/update dates ----
$ws.Range("U8").Value = "En este periodo no hubo personal contratado por honorarios."
$ws.Range("R8").Value = "Departamento de Recursos Humanos (UPP)"
$ws.Range("S8").Value = 44844   # 2022-10-10
$ws.Range("T8").Value = 44844   # 2022-10-10

# --- 6. Row 8 no longer needs the taller custom height ---------------------
$ws.Rows("8:8").AutoFit()

# --- 7. Remove now-stale hyperlinks (no more contract/legal URLs) ---------
$ws.Hyperlinks.Delete()

# --- 8. Extend the data-validation range for column D (catalog dropdown) --
$ws.Range("D8:D201").Validation.Delete()
$ws.Range("D8:D201").Validation.Add(3, 1, 1, "Hidden_13")
$ws.Range("D8:D201").Validation.ShowInput = $false

# --- 9. Drop the stale AutoFilter defined name -----------------------------
$wb.Names.Item("Reporte de Formatos!_FilterDatabase").Delete()

# --- 10. Re-fit columns whose content width changed ------------------------
$ws.Columns("D:D").AutoFit()
$ws.Columns("J:J").AutoFit()
$ws.Columns("M:M").AutoFit()
$ws.Columns("Q:Q").AutoFit()
$ws.Columns("S:S").AutoFit()
$ws.Columns("T:T").AutoFit()
$ws.Columns("U:U").AutoFit()
